$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New participant row appended by the Streamlit app (Juan, 2026-01-30 14:18:41)
$row = 24

$ws.Range("A$row").Value = "Juan_20260130_141841"
# B$row ("Grupo_Experimental") intentionally left blank for this participant
$ws.Range("C$row").Value = "Juan"
$ws.Range("D$row").Value = 23
$ws.Range("E$row").Value = "Male"
$ws.Range("F$row").Value = "2026-01-30 14:18:41"

$ws.Range("G$row").Value = @"
{
  "portion": 0.6,
  "diet": 0.8571428571428571,
  "salt": 0.4,
  "fat": 0.4,
  "natural": 0.6,
  "convenience": 0.4,
  "price": 0.2
}
"@

$ws.Range("H$row").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I$row").NumberFormat = "@"
$ws.Range("I$row").Value = "0.560"
$ws.Range("J$row").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Range("K$row").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("L$row").NumberFormat = "@"
$ws.Range("L$row").Value = "0.463"
$ws.Range("M$row").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

$ws.Range("N$row").Value = "Nongshim Shin Ramyun"
$ws.Range("O$row").NumberFormat = "@"
$ws.Range("O$row").Value = "0.430"
$ws.Range("P$row").Value = "Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio"

$ws.Range("Q$row").Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("R$row").NumberFormat = "@"
$ws.Range("R$row").Value = "0.718"
$ws.Range("S$row").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

$ws.Range("T$row").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("U$row").NumberFormat = "@"
$ws.Range("U$row").Value = "0.584"
$ws.Range("V$row").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Range("W$row").Value = "Annie’s Shells & White Cheddar"
$ws.Range("X$row").NumberFormat = "@"
$ws.Range("X$row").Value = "0.561"
$ws.Range("Y$row").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

$ws.Range("Z$row").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA$row").NumberFormat = "@"
$ws.Range("AA$row").Value = "0.780"
$ws.Range("AB$row").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

$ws.Range("AC$row").Value = "Kitchens of India Variety Pack"
$ws.Range("AD$row").NumberFormat = "@"
$ws.Range("AD$row").Value = "0.517"
$ws.Range("AE$row").Value = "Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad"

$ws.Range("AF$row").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AG$row").NumberFormat = "@"
$ws.Range("AG$row").Value = "0.504"
$ws.Range("AH$row").Value = "Portátil, saludable, fácil, buena textura, sabor suave"
